$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data refresh adds two new price records (row 58 and 59) for
# "Damasco" at the top of this block, pushing the previously existing
# records (old rows 58-103) down by two rows (new rows 60-105).
$ws.Rows.Item(58).Resize(2).Insert()

# New row 58: Damasco, Castle Brite, Especial - Provincia de Limarí
$ws.Cells.Item(58,1).Value  = 10
$ws.Cells.Item(58,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(58,3).Value  = "La Araucanía"
$ws.Cells.Item(58,4).Value  = "12/11/2023"
$ws.Cells.Item(58,5).Value  = 9
$ws.Cells.Item(58,6).Value  = "Fruta"
$ws.Cells.Item(58,7).Value  = 100103
$ws.Cells.Item(58,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(58,9).Value  = 100103003
$ws.Cells.Item(58,10).Value = "Damasco"
$ws.Cells.Item(58,11).Value = "Castle Brite"
$ws.Cells.Item(58,12).Value = "Especial"
$ws.Cells.Item(58,13).Value = 25
$ws.Cells.Item(58,14).Value = 15000
$ws.Cells.Item(58,15).Value = 15000
$ws.Cells.Item(58,16).Value = 15000
$ws.Cells.Item(58,17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(58,18).Value = "Provincia de Limarí"
$ws.Cells.Item(58,19).Value = 2143
$ws.Cells.Item(58,20).Value = 7

# New row 59: Damasco, Castle Brite, Primera - Región de O'Higgins
$ws.Cells.Item(59,1).Value  = 10
$ws.Cells.Item(59,2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(59,3).Value  = "La Araucanía"
$ws.Cells.Item(59,4).Value  = "12/11/2023"
$ws.Cells.Item(59,5).Value  = 9
$ws.Cells.Item(59,6).Value  = "Fruta"
$ws.Cells.Item(59,7).Value  = 100103
$ws.Cells.Item(59,8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(59,9).Value  = 100103003
$ws.Cells.Item(59,10).Value = "Damasco"
$ws.Cells.Item(59,11).Value = "Castle Brite"
$ws.Cells.Item(59,12).Value = "Primera"
$ws.Cells.Item(59,13).Value = 15
$ws.Cells.Item(59,14).Value = 30000
$ws.Cells.Item(59,15).Value = 30000
$ws.Cells.Item(59,16).Value = 30000
$ws.Cells.Item(59,17).Value = "$/caja 15 kilos"
$ws.Cells.Item(59,18).Value = "Región de O'Higgins"
$ws.Cells.Item(59,19).Value = 2000
$ws.Cells.Item(59,20).Value = 15
